$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Style tweak: the Phase/Task column style now right-aligns its
# values instead of centering them (it now holds numbers instead of
# the placeholder text "-"). Do this first, while the style is only
# referenced by the single existing task row, so the existing style
# record is updated in place rather than a new duplicate being made.
# ------------------------------------------------------------------
$ws.Range("F6").HorizontalAlignment = -4152

# ------------------------------------------------------------------
# Structural changes: the log table grows from 1 task entry to 3, and
# a "total hours" formula row is added above the table header.
# ------------------------------------------------------------------

# Step 1: Insert a new blank row at 5 (will hold the total-hours sum
# formula). This pushes the header row (old row 5) down to row 6 and
# the single task-data row (old row 6) down to row 7.
$ws.Rows("5:5").Insert()

# Step 2: Insert 2 new blank rows right after the task-data template
# row (row 7), to host the two additional task entries.
$ws.Rows("8:9").Insert()

# Step 3: Duplicate the formatting (number formats, fonts, fills,
# alignment...) of the template data row (row 7) into the two new
# rows 8 and 9, so all three task rows share identical formatting.
$ws.Range("A7:H7").Copy()
$ws.Range("A8:H8").PasteSpecial()
$ws.Range("A7:H7").Copy()
$ws.Range("A9:H9").PasteSpecial()

# Step 4: Remove the old trailing blank row + old sum-formula row
# (now shifted down to 10:11), since the sum formula moves to row 5.
$ws.Rows("10:11").Delete()

# Fix up row heights: the new total-hours row uses the sheet's normal
# row height, while the task rows keep the slightly taller row height
# used by the header/table rows.
$ws.Rows("5:5").RowHeight = 13.75
$ws.Rows("7:9").RowHeight = 14.15

# ------------------------------------------------------------------
# Content changes
# ------------------------------------------------------------------

# Header date field (top of sheet) updated.
$ws.Range("F1").Value2 = 41902

# Row 5: total hours formula (sum of the delta-time column across the
# three task rows, converted from minutes to hours).
$ws.Range("E5").Formula = "=SUM(E7:E9)/60"

# Row 7: first task entry (new).
$ws.Range("A7").Value2 = 41902
$ws.Range("B7").Value2 = 0.631944444444444
$ws.Range("C7").Value2 = 0.684027777777778
$ws.Range("D7").Value2 = 14
$ws.Range("E7").Formula = "=((HOUR(C7)-HOUR(B7))*60)+(MINUTE(C7)-MINUTE(B7))-D7"
$ws.Range("F7").Value2 = 2
$ws.Range("G7").ClearContents()
$ws.Range("H7").Value2 = "Realizar el lanzamiento del ciclo #1 de TSPi."

# Row 8: second task entry (new).
$ws.Range("A8").Value2 = 41902
$ws.Range("B8").Value2 = 0.6875
$ws.Range("C8").Value2 = 0.725694444444444
$ws.Range("D8").Value2 = 10
$ws.Range("E8").Formula = "=((HOUR(C8)-HOUR(B8))*60)+(MINUTE(C8)-MINUTE(B8))-D8"
$ws.Range("F8").Value2 = 3
$ws.Range("G8").ClearContents()
$ws.Range("H8").Value2 = "Definir la estrategía de desarrollo del ciclo #1 de TSPi."

# Row 9: the pre-existing task entry (unchanged date/time/comment),
# but its Phase/Task column switches from the placeholder text "-" to
# the numeric value 1, and its formula is re-applied (pasting above
# flattened it to a static value).
$ws.Range("A9").Value2 = 41904
$ws.Range("B9").Value2 = 0.340277777777778
$ws.Range("C9").Value2 = 0.354166666666667
$ws.Range("D9").Value2 = 0
$ws.Range("E9").Formula = "=((HOUR(C9)-HOUR(B9))*60)+(MINUTE(C9)-MINUTE(B9))-D9"
$ws.Range("F9").Value2 = 1
$ws.Range("G9").ClearContents()
$ws.Range("H9").Value2 = "Video tutorial de Github."

# Selection moves to E6 (the header's "Delta Time" cell).
$ws.Range("E6").Select()
